$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 495, shifting existing rows 495:598 down to 496:599
$ws.Rows.Item(495).Insert()

# Populate the newly inserted row 495 with the new data record
$ws.Cells.Item(495, 1).Value = 8
$ws.Cells.Item(495, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(495, 3).Value = "Coquimbo"
$ws.Cells.Item(495, 4).Value = 45211
$ws.Cells.Item(495, 5).Value = 4
$ws.Cells.Item(495, 6).Value = 100112032
$ws.Cells.Item(495, 7).Value = "Zapallo italiano"
$ws.Cells.Item(495, 8).Value = "Sin especificar"
$ws.Cells.Item(495, 9).Value = "Primera"
$ws.Cells.Item(495, 10).Value = 300
$ws.Cells.Item(495, 11).Value = 13000
$ws.Cells.Item(495, 12).Value = 14000
$ws.Cells.Item(495, 13).Value = 13500
$ws.Cells.Item(495, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(495, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(495, 16).Value = 270
$ws.Cells.Item(495, 17).Value = 50
$ws.Cells.Item(495, 18).Value = "Hortaliza"
